$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "For testing" label in E3 is replaced with a new label describing
# the dataset split used for summarization evaluation.
$ws.Range("E3").Value = "For summarization evaluation"

# Column E is widened to fit the new (longer) label.
$ws.Columns("E").AutoFit() | Out-Null

# Leave the selection where the editor last left it when saving the file.
$ws.Range("E13").Select() | Out-Null
